$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)  # Classi
$ws2 = $wb.Worksheets.Item(2)  # Alunni in presenza
$ws3 = $wb.Worksheets.Item(3)  # Alunni
$ws4 = $wb.Worksheets.Item(4)  # Personale scolastico

$newLabel = "7  - 12 marzo 2022"

# ---------------------------------------------------------------------
# Sheet 1 "Classi": add row 10, copying the per-column formatting from
# row 9 (the previous last row) so styles match, then writing values.
# ---------------------------------------------------------------------
foreach ($col in @("A","B","C","D","E","F","G","H","I","J","K","N","O")) {
    $ws1.Range("${col}9").Copy()
    $ws1.Range("${col}10").PasteSpecial(-4122)
}
$ws1.Range("A10").Value = $newLabel
$ws1.Range("B10").Value = 5827
$ws1.Range("C10").Value = 8157
$ws1.Range("D10").Value = 0.71400000000000008
$ws1.Range("E10").Value = 376491
$ws1.Range("F10").Value = 268352
$ws1.Range("G10").Value = 0.71299999999999997
$ws1.Range("H10").Value = 268076
$ws1.Range("I10").Value = 20809
$ws1.Range("J10").Value = 0.99900000000000011
$ws1.Range("K10").Value = 0.078
$ws1.Range("N10").Value = 276
$ws1.Range("O10").Value = 0.001

# ---------------------------------------------------------------------
# Sheet 2 "Alunni in presenza": add row 10, based on row 9 formatting.
# ---------------------------------------------------------------------
foreach ($col in @("A","B","C","D","E","F")) {
    $ws2.Range("${col}9").Copy()
    $ws2.Range("${col}10").PasteSpecial(-4122)
}
$ws2.Range("A10").Value = $newLabel
$ws2.Range("B10").Value = 7393374
$ws2.Range("C10").Value = 5248017
$ws2.Range("D10").Value = 0.71
$ws2.Range("E10").Value = 5098082
$ws2.Range("F10").Value = 0.97099999999999997

# ---------------------------------------------------------------------
# Sheet 3 "Alunni": add rows 34, 35, 36 (row 33 stays empty, matching
# the existing gap pattern between date blocks), based on rows 30/31/32.
# ---------------------------------------------------------------------
foreach ($col in @("A","B","C","D","E","F")) {
    $ws3.Range("${col}30").Copy()
    $ws3.Range("${col}34").PasteSpecial(-4122)
}
$ws3.Range("A34").Value = $newLabel
$ws3.Range("B34").Value = "Infanzia"
$ws3.Range("C34").Value = 603087
$ws3.Range("D34").Value = 590099
$ws3.Range("E34").Value = 12988
$ws3.Range("F34").Value = 0.022000000000000002

foreach ($col in @("A","B","C","D","E","F")) {
    $ws3.Range("${col}31").Copy()
    $ws3.Range("${col}35").PasteSpecial(-4122)
}
$ws3.Range("A35").Value = $newLabel
$ws3.Range("B35").Value = "Primaria"
$ws3.Range("C35").Value = 1669122
$ws3.Range("D35").Value = 1622263
$ws3.Range("E35").Value = 46859
$ws3.Range("F35").Value = 0.027999999999999997

foreach ($col in @("A","B","C","D","E","F")) {
    $ws3.Range("${col}32").Copy()
    $ws3.Range("${col}36").PasteSpecial(-4122)
}
$ws3.Range("A36").Value = $newLabel
$ws3.Range("B36").Value = "Sec. 1° e 2° Grado"
$ws3.Range("C36").Value = 2975808
$ws3.Range("D36").Value = 2885720
$ws3.Range("E36").Value = 90088
$ws3.Range("F36").Value = 0.03

# ---------------------------------------------------------------------
# Sheet 4 "Personale scolastico": add row 10, based on row 9 formatting.
# ---------------------------------------------------------------------
foreach ($col in @("A","B","C","D","E","F","G","H","I","J","K")) {
    $ws4.Range("${col}9").Copy()
    $ws4.Range("${col}10").PasteSpecial(-4122)
}
$ws4.Range("A10").Value = $newLabel
$ws4.Range("B10").Value = 775867
$ws4.Range("C10").Value = 548032
$ws4.Range("D10").Value = 0.70599999999999996
$ws4.Range("E10").Value = 530432
$ws4.Range("F10").Value = 0.96799999999999997
$ws4.Range("G10").Value = 204526
$ws4.Range("H10").Value = 144880
$ws4.Range("I10").Value = 0.70799999999999996
$ws4.Range("J10").Value = 141428
$ws4.Range("K10").Value = 0.97599999999999998

# ---------------------------------------------------------------------
# Update selections on each sheet to match the new last-cell positions.
# ---------------------------------------------------------------------
[void]$ws1.Range("A10").Select()
[void]$ws2.Range("A10").Select()
[void]$ws4.Range("A10").Select()

# Activate "Alunni" (3rd tab) last so it becomes the workbook's active
# tab/sheet and picks up tabSelected, then set its selection too.
[void]$ws3.Activate()
[void]$ws3.Range("A36").Select()
